$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from SCD0336 -> SCD0025
$ws.Name = "SCD0025"

# Update the TC_ID cell (B2) from the old "DGS-351" placeholder to the
# new test-case id "SCD0025-006"
$ws.Range("B2").Value = "SCD0025-006"

# Widen column B so the longer TC_ID text fits (Excel auto-fit would
# grow the "best fit" column from 9 chars to ~12.43 chars)
$ws.Columns.Item(2).ColumnWidth = 11.665

# Move the active selection from B2 to B3, as captured in the saved view
$ws.Range("B3").Select()
